$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 210; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value2
    if ($cur -eq 45180) {
        $cell.Value2 = 45181
    }
}
